$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.076453371673271
$ws.Range("D2").Value = 1.081735160211898
$ws.Range("E2").Value = 1.08047521814624
$ws.Range("F2").Value = 1.091889822404316
$ws.Range("I2").Value = 1.044389096754767
$ws.Range("J2").Value = 1.081352752234992
$ws.Range("K2").Value = 1.08440520869858
$ws.Range("L2").Value = 1.083148555641206
$ws.Range("M2").Value = 1.094533661629578
$ws.Range("N2").Value = 1.08288839717606

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.077979213603912
$ws.Range("D3").Value = 1.083133123019333
$ws.Range("E3").Value = 1.08184898137148
$ws.Range("F3").Value = 1.093381887774578
$ws.Range("I3").Value = 1.044668929334359
$ws.Range("J3").Value = 1.082535748191552
$ws.Range("K3").Value = 1.085620270200645
$ws.Range("L3").Value = 1.084339236193553
$ws.Range("M3").Value = 1.095844513708981
$ws.Range("N3").Value = 1.08407307312257

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.078965298278686
$ws.Range("D4").Value = 1.084036784077385
$ws.Range("E4").Value = 1.082737039370468
$ws.Range("F4").Value = 1.094346636689477
$ws.Range("I4").Value = 1.044848010318628
$ws.Range("J4").Value = 1.08329954189866
$ws.Range("K4").Value = 1.086405045208015
$ws.Range("L4").Value = 1.085108286210864
$ws.Range("M4").Value = 1.096691483502999
$ws.Range("N4").Value = 1.08483795150435

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.079379559127611
$ws.Range("D5").Value = 1.084416469972735
$ws.Range("E5").Value = 1.0831101789063
$ws.Range("F5").Value = 1.09475205149247
$ws.Range("I5").Value = 1.044922821161198
$ws.Range("J5").Value = 1.083620242894124
$ws.Range("K5").Value = 1.086734622475242
$ws.Range("L5").Value = 1.085431264792939
$ws.Range("M5").Value = 1.097047258481703
$ws.Range("N5").Value = 1.085159107932012

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.07944909859432
$ws.Range("D6").Value = 1.084480208566986
$ws.Range("E6").Value = 1.083172819105278
$ws.Range("F6").Value = 1.094820112883686
$ws.Range("I6").Value = 1.04493535443052
$ws.Range("J6").Value = 1.083674066821799
$ws.Range("K6").Value = 1.086789940036601
$ws.Range("L6").Value = 1.085485475111666
$ws.Range("M6").Value = 1.097106977754795
$ws.Range("N6").Value = 1.085213008295839

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.078970834779903
$ws.Range("D7").Value = 1.084041858289653
$ws.Range("E7").Value = 1.082742026058429
$ws.Range("F7").Value = 1.094352054505289
$ws.Range("I7").Value = 1.044849011808618
$ws.Range("J7").Value = 1.083303828678004
$ws.Range("K7").Value = 1.086409450373524
$ws.Range("L7").Value = 1.08511260315704
$ws.Range("M7").Value = 1.096696238522751
$ws.Range("N7").Value = 1.084842244371413

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.076969297474734
$ws.Range("D8").Value = 1.082207801042638
$ws.Range("E8").Value = 1.08093966898433
$ws.Range("F8").Value = 1.092394225021467
$ws.Range("I8").Value = 1.044484080289923
$ws.Range("J8").Value = 1.081752903046705
$ws.Range("K8").Value = 1.08481614802708
$ws.Range("L8").Value = 1.083551244647722
$ws.Range("M8").Value = 1.094976930141865
$ws.Range("N8").Value = 1.083289116247823

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.073432543446995
$ws.Range("D9").Value = 1.078968699347146
$ws.Range("E9").Value = 1.077756860977216
$ws.Range("F9").Value = 1.08893851503276
$ws.Range("I9").Value = 1.043825721728621
$ws.Range("J9").Value = 1.079006843356728
$ws.Range("K9").Value = 1.081997197169716
$ws.Range("L9").Value = 1.080788980513003
$ws.Range("M9").Value = 1.091937549869803
$ws.Range("N9").Value = 1.08053915683811

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.071067683117708
$ws.Range("D10").Value = 1.076804058358022
$ws.Range("E10").Value = 1.075630051257584
$ws.Range("F10").Value = 1.086630450551555
$ws.Range("I10").Value = 1.043376433671809
$ws.Range("J10").Value = 1.077166978149484
$ws.Range("K10").Value = 1.080109935228008
$ws.Range("L10").Value = 1.078939785638344
$ws.Range("M10").Value = 1.089904398202298
$ws.Range("N10").Value = 1.078696678811236

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.070041907404142
$ws.Range("D11").Value = 1.075865421383322
$ws.Range("E11").Value = 1.074707869620288
$ws.Range("F11").Value = 1.085629937530829
$ws.Range("I11").Value = 1.04317940315948
$ws.Range("J11").Value = 1.076368049675652
$ws.Range("K11").Value = 1.07929076728815
$ws.Range("L11").Value = 1.078137168525657
$ws.Range("M11").Value = 1.089022306388657
$ws.Range("N11").Value = 1.077896615767333

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.069660613259371
$ws.Range("D12").Value = 1.075516562490586
$ws.Range("E12").Value = 1.074365134313037
$ws.Range("F12").Value = 1.085258128834553
$ws.Range("I12").Value = 1.043105841891863
$ws.Range("J12").Value = 1.076070946843195
$ws.Range("K12").Value = 1.078986189218921
$ws.Range("L12").Value = 1.077838748484427
$ws.Range("M12").Value = 1.088694392390433
$ws.Range("N12").Value = 1.077599091014775

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.069742414747821
$ws.Range("D13").Value = 1.075591403402874
$ws.Range("E13").Value = 1.074438661179036
$ws.Range("F13").Value = 1.085337891038864
$ws.Range("I13").Value = 1.043121638050045
$ws.Range("J13").Value = 1.076134692136046
$ws.Range("K13").Value = 1.07905153604777
$ws.Range("L13").Value = 1.077902773925943
$ws.Range("M13").Value = 1.088764743195832
$ws.Range("N13").Value = 1.077662926833253

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.070010395178722
$ws.Range("D14").Value = 1.075836588868738
$ws.Range("E14").Value = 1.074679543050683
$ws.Range("F14").Value = 1.085599207293046
$ws.Range("I14").Value = 1.043173330230462
$ws.Range("J14").Value = 1.076343498142994
$ws.Range("K14").Value = 1.07926559698108
$ws.Range("L14").Value = 1.078112507031731
$ws.Range("M14").Value = 1.088995206377814
$ws.Range("N14").Value = 1.077872029368683

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.070175470018578
$ws.Range("D15").Value = 1.075987628042058
$ws.Range("E15").Value = 1.074827932146318
$ws.Range("F15").Value = 1.085760189632086
$ws.Range("I15").Value = 1.043205129692683
$ws.Range("J15").Value = 1.076472104474944
$ws.Range("K15").Value = 1.07939744669147
$ws.Range("L15").Value = 1.078241691596461
$ws.Range("M15").Value = 1.089137166920404
$ws.Range("N15").Value = 1.078000818336375

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.071135722330393
$ws.Range("D16").Value = 1.076866323915393
$ws.Range("E16").Value = 1.075691226307712
$ws.Range("F16").Value = 1.086696827370027
$ws.Range("I16").Value = 1.043389457403169
$ws.Range("J16").Value = 1.077219952345009
$ws.Range("K16").Value = 1.080164258625149
$ws.Range("L16").Value = 1.078993012048757
$ws.Range("M16").Value = 1.089962902830151
$ws.Range("N16").Value = 1.078749728236194

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.071737581801166
$ws.Range("D17").Value = 1.07741714474913
$ws.Range("E17").Value = 1.076232406047211
$ws.Range("F17").Value = 1.087284054148994
$ws.Range("I17").Value = 1.043504414447964
$ws.Range("J17").Value = 1.077688449347878
$ws.Range("K17").Value = 1.080644727213169
$ws.Range("L17").Value = 1.079463781438852
$ws.Range("M17").Value = 1.09048039883769
$ws.Range("N17").Value = 1.079218890558545

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.072088465514913
$ws.Range("D18").Value = 1.077738301072233
$ws.Range("E18").Value = 1.076547946023755
$ws.Range("F18").Value = 1.087626467401697
$ws.Range("I18").Value = 1.043571227231088
$ws.Range("J18").Value = 1.077961498618653
$ws.Range("K18").Value = 1.08092478667418
$ws.Range("L18").Value = 1.079738190026676
$ws.Range("M18").Value = 1.090782079655304
$ws.Range("N18").Value = 1.079492327590604

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.072208079093679
$ws.Range("D19").Value = 1.077847785554328
$ws.Range("E19").Value = 1.076655516681935
$ws.Range("F19").Value = 1.087743203684831
$ws.Range("I19").Value = 1.043593968060493
$ws.Range("J19").Value = 1.078054564734129
$ws.Range("K19").Value = 1.081020247725631
$ws.Range("L19").Value = 1.079831725374699
$ws.Range("M19").Value = 1.090884917093504
$ws.Range("N19").Value = 1.079585525870638

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.071673025698039
$ws.Range("D20").Value = 1.077358060192816
$ws.Range("E20").Value = 1.076174355149028
$ws.Range("F20").Value = 1.087221061311231
$ws.Range("I20").Value = 1.043492105453742
$ws.Range("J20").Value = 1.077638206556462
$ws.Range("K20").Value = 1.080593197153329
$ws.Range("L20").Value = 1.079413291337115
$ws.Range("M20").Value = 1.0904248936178
$ws.Range("N20").Value = 1.079168576416602

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.069931489292127
$ws.Range("D21").Value = 1.075764393640165
$ws.Range("E21").Value = 1.074608614808672
$ws.Range("F21").Value = 1.08552226102213
$ws.Range("I21").Value = 1.043158118548878
$ws.Range("J21").Value = 1.076282019560499
$ws.Range("K21").Value = 1.079202569777805
$ws.Range("L21").Value = 1.07805075397973
$ws.Range("M21").Value = 1.088927348108185
$ws.Range("N21").Value = 1.077810463479548

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.068834918705558
$ws.Range("D22").Value = 1.074761189169079
$ws.Range("E22").Value = 1.073623034368579
$ws.Range("F22").Value = 1.084453149647899
$ws.Range("I22").Value = 1.042945955042688
$ws.Range("J22").Value = 1.075427330358156
$ws.Range("K22").Value = 1.07832647301919
$ws.Range("L22").Value = 1.077192378286995
$ws.Range("M22").Value = 1.087984240404838
$ws.Range("N22").Value = 1.076954560520504

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.069416385837993
$ws.Range("D23").Value = 1.075293123238858
$ws.Range("E23").Value = 1.074145619309949
$ws.Range("F23").Value = 1.085020003888733
$ws.Range("I23").Value = 1.043058633513966
$ws.Range("J23").Value = 1.07588060908637
$ws.Range("K23").Value = 1.078791076938322
$ws.Range("L23").Value = 1.077647582213408
$ws.Range("M23").Value = 1.088484348050315
$ws.Range("N23").Value = 1.077408482956504

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.071702196350085
$ws.Range("D24").Value = 1.077384758356831
$ws.Range("E24").Value = 1.076200586225363
$ws.Range("F24").Value = 1.087249525391569
$ws.Range("I24").Value = 1.043497668099208
$ws.Range("J24").Value = 1.077660909782326
$ws.Range("K24").Value = 1.080616481955829
$ws.Range("L24").Value = 1.07943610620721
$ws.Range("M24").Value = 1.09044997455022
$ws.Range("N24").Value = 1.079191311883651

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.074348084619601
$ws.Range("D25").Value = 1.079806981819765
$ws.Range("E25").Value = 1.078580537692518
$ws.Range("F25").Value = 1.089832621189307
$ws.Range("I25").Value = 1.043997746234033
$ws.Range("J25").Value = 1.079718355531191
$ws.Range("K25").Value = 1.082727342357795
$ws.Range("L25").Value = 1.081504422825448
$ws.Range("M25").Value = 1.092724492561917
$ws.Range("N25").Value = 1.081251679441472

